$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: add "I0" in I1 and "IF" in J1, matching the styling of the
#     existing header cells (e.g. H1) ---
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data: fill I2:J76 with the per-row (I, J) values ---
$data = @(
    @(8,8),
    @(7,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(10,10),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(9,10),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,10),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(4,4),
    @(3,3)
)

for ($k = 0; $k -lt $data.Length; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $data[$k][0]
    $ws.Cells.Item($r, 10).Value = $data[$k][1]
}
